$d = $word.ActiveDocument

# 1. Table: add row/table justification w:jc val="start" to tblPr.
$t = $d.Tables(1)
$t.Alignment = "start"

# 2. New paragraph style "AbstractTitle" (display name "Abstract Title"),
#    inserted before the "Abstract" style definition.
$absTitle = $d.Styles.Add("AbstractTitle", 1)
$absTitle.NameLocal = "Abstract Title"
$absTitle.BaseStyle = "Normal"
$absTitle.NextParagraphStyle = "Abstract"
$absTitle.QuickStyle = $true
$absTitle.ParagraphFormat.KeepWithNext = $true
$absTitle.ParagraphFormat.KeepTogether = $true
$absTitle.ParagraphFormat.Alignment = 1
$absTitle.ParagraphFormat.SpaceBefore = 15
$absTitle.ParagraphFormat.SpaceAfter = 0
$absTitle.Font.Size = 10
$absTitle.Font.SizeBi = 10
$absTitle.Font.Bold = $true
$absTitle.Font.Color = 9067060

# 3. "Abstract" style: reduce space-before from 300 (15pt) to 100 (5pt) twips.
$abs = $d.Styles("Abstract")
$abs.ParagraphFormat.SpaceBefore = 5

# 4. "ImportTok" character style: add green + bold.
$importTok = $d.Styles("ImportTok")
$importTok.Font.Color = 32768
$importTok.Font.Bold = $true

# 5. "BuiltInTok" character style: add green.
$builtInTok = $d.Styles("BuiltInTok")
$builtInTok.Font.Color = 32768

Write-Output "done"
